$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: paste formatting (styles) for new rows 13-18 from matching template rows ---
$ws.Range("B6:E6").Copy()
$ws.Range("B13:E13").PasteSpecial(-4122)
$ws.Range("A4:E4").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)
$ws.Range("B6:E6").Copy()
$ws.Range("B15:E15").PasteSpecial(-4122)
$ws.Range("A4:E4").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)
$ws.Range("B6:E6").Copy()
$ws.Range("B17:E17").PasteSpecial(-4122)
$ws.Range("B6:E6").Copy()
$ws.Range("B18:E18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Step 2: set B column (row numbers) ---
$ws.Range("B13").Value = 106
$ws.Range("B14").Value = 109
$ws.Range("B15").Value = 51
$ws.Range("B16").Value = 54
$ws.Range("B17").Value = 24
$ws.Range("B18").Value = 27

# --- Step 3: set C/D/E text values in the exact order to match shared-string index assignment ---
$ws.Range("C13").Value = " T-Team [CS:X]Ch-Charm[CR] is here,\nright?! I mean, at the guild?!"
$ws.Range("C14").Value = " Whoaaah! Whoooaaah![K]\nWhoo-hoo-whoooaaah!"
$ws.Range("D13").Value = " Здесь К-Команда [CS:X]Ш[CR]-[CS:X]Шарм[CR], да?!\nВ гильдии?!"
$ws.Range("D14").Value = " Вааау! Ваааааау![K]\nВаау-вау-ваааааау!"
$ws.Range("E13").Value = " Èäåòû Ë-Ëïíàîäà [CS:X]Š[CR]-[CS:X]Šàñí[CR], äà?!\nÂ ãéìûäéé?!"
$ws.Range("E14").Value = " Âàààô! Âààààààô![K]\nÂààô-âàô-âààààààô!"
$ws.Range("C15").Value = " You throw everything you have\ninto this, you two!"
$ws.Range("C16").Value = " You set that [CS:N]Darkrai[CR] straight!"
$ws.Range("C17").Value = " All right! You two are heroes!"
$ws.Range("C18").Value = " It\'ll get peaceful in Treasure\nTown too. Thanks for that!"
$ws.Range("D15").Value = " Покажите, на что вы, ребята,\nспособны!"
$ws.Range("D16").Value = " Укажите [CS:N]Даркраю[CR] его место!"
$ws.Range("D17").Value = " Чудесно! Вы настоящие герои!"
$ws.Range("D18").Value = " Скоро в Город Сокровищ придёт\nпокой. Спасибо вам за всё!"
$ws.Range("E15").Value = " Ðïëàçéóå, îà œóï âú, ñåáÿóà,\nòðïòïáîú!"
$ws.Range("E16").Value = " Ôëàçéóå [CS:N]Äàñëñàý[CR] åãï íåòóï!"
$ws.Range("E17").Value = " Œôäåòîï! Âú îàòóïÿþéå ãåñïé!"
$ws.Range("E18").Value = " Òëïñï â Ãïñïä Òïëñïâéþ ðñéäæó\nðïëïê. Òðàòéáï âàí èà âòæ!"

# --- Step 4: row heights ---
$ws.Rows.Item(13).RowHeight = 21.6
$ws.Rows.Item(14).RowHeight = 21.6
$ws.Rows.Item(15).RowHeight = 21.6
$ws.Rows.Item(16).RowHeight = 24
$ws.Rows.Item(18).RowHeight = 21.6

# --- Step 5: selection ---
$ws.Range("C18").Select()

